# Automatic update of files.
#
# 1) Column C ("Förändrad") bumps from 45184 -> 45186 for every data row.
# 2) Every HYPERLINK(...) formula in columns S:Y gains a second argument
#    (the friendly link text), equal to the row's designation in column A,
#    e.g. HYPERLINK("...A 11286-2021.xlsx") -> HYPERLINK("...A 11286-2021.xlsx", "A 11286-2021")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$firstDataRow = 2
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

# 1) Bump the "Förändrad" date column (C) for every data row in one shot.
$ws.Range("C$firstDataRow`:C$lastRow").Value = 45186

# 2) Add the designation as the 2nd HYPERLINK() argument, columns S (19) .. Y (25).
$firstCol = 19
$lastCol = 25

for ($row = $firstDataRow; $row -le $lastRow; $row++) {
    $designation = $ws.Cells.Item($row, 1).Value()
    for ($col = $firstCol; $col -le $lastCol; $col++) {
        $cell = $ws.Cells.Item($row, $col)
        if ($cell.HasFormula()) {
            $f = $cell.Formula()
            if ($f.Contains("HYPERLINK") -and -not $f.Contains(",")) {
                $newFormula = $f.Replace('")', '", "' + $designation + '")')
                $cell.Formula = $newFormula
            }
        }
    }
}
